$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-09-30 00:00:00"
$ws.Range("O2").Value = 70632634.43000001
$ws.Range("P2").Value = 189787133.01
$ws.Range("Q2").Value = 116660667.33
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = 70373515.78
$ws.Range("T2").Value = 70373515.78
$ws.Range("U2").Value = ""
$ws.Range("V2").Value = 9891868.42
$ws.Range("W2").Value = 13918491.52
$ws.Range("X2").Value = 1568600.86
$ws.Range("Y2").Value = 82564264.17
$ws.Range("Z2").Value = 82434264.17
$ws.Range("AA2").Value = 11801629.74
$ws.Range("AG2").Value = 2749156.13
$ws.Range("AP2").Value = ""
$ws.Range("AQ2").Value = ""
$ws.Range("AR2").Value = ""
$ws.Range("AS2").Value = 58970430.5
$ws.Range("AT2").Value = ""
